# Adding overage learner tab. Change gender labels.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "overage" sheet as a copy of "access" (keeps column
#    widths / fonts / fills identical), placed right after "access".
# ---------------------------------------------------------------------------
$access = $wb.Worksheets.Item("access")
$access.Copy([System.Reflection.Missing]::Value, $access) | Out-Null
$overage = $wb.Worksheets.Item(2)
$overage.Name = "overage"

# ---------------------------------------------------------------------------
# 2. Re-purpose the copied cells for the overage-learners content.
# ---------------------------------------------------------------------------
$overage.Range("A2").Value = "overage"
$overage.Range("B2").Value = "Analysis of overage learners "
$overage.Range("F2").Value = "Overage learners"

# G2 / G3 reuse the "% ... at least 2 years above the intended age ..."
# descriptions that used to live on level1/level2 G3 - copy their format
# (no-fill, style index 3) along with the text itself.
$level1 = $wb.Worksheets.Item("level1")
$level2 = $wb.Worksheets.Item("level2")

$level1.Range("G3").Copy() | Out-Null
$overage.Range("G2").PasteSpecial(-4122) | Out-Null

$level2.Range("G3").Copy() | Out-Null
$overage.Range("G3").PasteSpecial(-4122) | Out-Null

# G4: blank cell with a new, slightly lighter fill (white/FFEFC1).
$overage.Range("G4").ClearContents() | Out-Null
$overage.Range("G4").Font.Name = "Segoe UI"
$overage.Range("G4").Font.Size = 10
$overage.Range("G4").Interior.Color = 16777215
$overage.Range("G4").VerticalAlignment = -4108

$overage.Range("C3").Value = "Girls"
$overage.Range("C4").Value = "Boys"

# Now that level1/level2's G3 text has been relocated onto "overage",
# clear it there (formatting/style is preserved, only the text goes away).
$level1.Range("G3").ClearContents() | Out-Null
$level2.Range("G3").ClearContents() | Out-Null

$overage.Range("C:C").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Gender-label rename ("Female / woman" -> "Girls", "Male / man" -> "Boys")
#    on every remaining sheet, plus selecting column C (matches the new
#    selection state recorded for every tab in the workbook).
# ---------------------------------------------------------------------------
$access.Range("C3").Value = "Girls"
$access.Range("C4").Value = "Boys"
$access.Range("C:C").Select() | Out-Null

$outOfSchool = $wb.Worksheets.Item("out_of_school")
$outOfSchool.Range("C3").Value = "Girls"
$outOfSchool.Range("C4").Value = "Boys"
$outOfSchool.Range("C:C").Select() | Out-Null

$ece = $wb.Worksheets.Item("ece")
$ece.Range("C3").Value = "Girls"
$ece.Range("C4").Value = "Boys"
$ece.Range("C:C").Select() | Out-Null

$level1.Range("C3").Value = "Girls"
$level1.Range("C4").Value = "Boys"
$level1.Range("C:C").Select() | Out-Null

$level2.Range("C3").Value = "Girls"
$level2.Range("C4").Value = "Boys"
$level2.Range("C:C").Select() | Out-Null

$level3 = $wb.Worksheets.Item("level3")
$level3.Range("C3").Value = "Girls"
$level3.Range("C4").Value = "Boys"
$level3.Range("C:C").Select() | Out-Null

# level4 is the tab that ends up active/selected, so it is handled last.
$level4 = $wb.Worksheets.Item("level4")
$level4.Range("C3").Value = "Girls"
$level4.Range("C4").Value = "Boys"
$level4.Range("H17").Select() | Out-Null
